$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4-12 (they are removed entirely in the new version)
$ws.Range("A4:E12").EntireRow.Delete()

# Update row 2 values
$ws.Range("B2").Value = "fregrwegrew"
$ws.Range("C2").Value = "суммы"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1"
$ws.Range("D2").Style = "Normal"

# Update row 3 values
$ws.Range("B3").Value = "qqq"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0"
$ws.Range("D3").Style = "Normal"

# Add new column F header and data, matching style of existing header cells (s="1")
$ws.Range("F1").Value = "Дата"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Value = "26.04.2021"
$ws.Range("F3").Value = "26.04.2021"
